$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# A3 value changes from "Demo1234" to "Demo12" (matching A2/A4)
$ws.Range("A3").Value = "Demo12"

# Update the selection to match the diff (activeCell B8, sqref B8)
$ws.Range("B8").Select()
